$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to be stored as text
# (prevents Excel auto-converting numeric-looking strings like "321.36"
# or "5.16%" into numbers/percentages), then reset the style so no
# quotePrefix formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '321.36'
Set-TextValue $ws.Range("E2") '5.16%'
Set-TextValue $ws.Range("D3") '36.19'
Set-TextValue $ws.Range("E3") '-0.14%'
Set-TextValue $ws.Range("D4") '5.125'
Set-TextValue $ws.Range("E4") '1.31%'
Set-TextValue $ws.Range("D5") '0.08144'
Set-TextValue $ws.Range("E5") '3.65%'
Set-TextValue $ws.Range("D6") '2.151'
Set-TextValue $ws.Range("E6") '-1.07%'
Set-TextValue $ws.Range("D7") '8.041'
Set-TextValue $ws.Range("E7") '1.54%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D8") '0.9277'
Set-TextValue $ws.Range("E8") '0.98%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D9") '0.1011'
Set-TextValue $ws.Range("E9") '3.96%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D10") '0.1885'
Set-TextValue $ws.Range("E10") '1.16%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D11") '0.09158'
Set-TextValue $ws.Range("E11") '5.40%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D12") '0.03590'
Set-TextValue $ws.Range("E12") '2.96%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D13") '0.09924'
Set-TextValue $ws.Range("E13") '0.00%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D14") '0.001431'
Set-TextValue $ws.Range("E14") '-0.42%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D15") '0.005650'
Set-TextValue $ws.Range("E15") '-0.08%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D16") '3.451'
Set-TextValue $ws.Range("E16") '-0.29%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D17") '4.142'
Set-TextValue $ws.Range("E17") '1.16%'
Set-TextValue $ws.Range("E18") '15.50%'
Set-TextValue $ws.Range("D19") '0.3373'
Set-TextValue $ws.Range("E19") '-1.53%'
Set-TextValue $ws.Range("D20") '0.1311'
Set-TextValue $ws.Range("E20") '-3.68%'
Set-TextValue $ws.Range("D21") '5.061'
Set-TextValue $ws.Range("E21") '5.89%'
Set-TextValue $ws.Range("D22") '0.2190'
Set-TextValue $ws.Range("E22") '-0.74%'
Set-TextValue $ws.Range("D23") '0.04609'
Set-TextValue $ws.Range("E23") '1.01%'
Set-TextValue $ws.Range("D24") '0.001243'
Set-TextValue $ws.Range("E24") '0.71%'
Set-TextValue $ws.Range("D25") '0.004730'
Set-TextValue $ws.Range("E25") '-7.21%'
Set-TextValue $ws.Range("D26") '0.0001301'
Set-TextValue $ws.Range("E26") '-7.12%'
Set-TextValue $ws.Range("D27") '0.0004504'
Set-TextValue $ws.Range("E27") '-5.20%'
Set-TextValue $ws.Range("D39") '0.02019'
Set-TextValue $ws.Range("E39") '10.36%'
Set-TextValue $ws.Range("E40") '5.26%'
Set-TextValue $ws.Range("D41") '0.007811'
Set-TextValue $ws.Range("E41") '0.40%'
Set-TextValue $ws.Range("D42") '0.1402'
Set-TextValue $ws.Range("E42") '0.34%'
Set-TextValue $ws.Range("D43") '0.007808'
Set-TextValue $ws.Range("E43") '0.97%'
Set-TextValue $ws.Range("E44") '-6.98%'
Set-TextValue $ws.Range("D45") '0.01218'
Set-TextValue $ws.Range("E45") '9.82%'
Set-TextValue $ws.Range("D46") '0.00006484'
Set-TextValue $ws.Range("E46") '2.89%'
Set-TextValue $ws.Range("E47") '0.06%'
Set-TextValue $ws.Range("E48") '20.49%'
Set-TextValue $ws.Range("E49") '-4.93%'
Set-TextValue $ws.Range("D50") '0.00002102'
Set-TextValue $ws.Range("E50") '0.06%'
Set-TextValue $ws.Range("E51") '0.06%'
